# LOQ4267.xlsx content update (site build 2023-04-12):
#  - Portuguese "Objetivos" text added, "Docentes responsaveis" value moved to its own row,
#    "Programa resumido"/"Programa"/"Metodo"/"Bibliografia" values replaced with real course text
#    (previously placeholder/duplicated values), and a new "Bibliografia" text row appended.
#  - Rows 10-21 reshuffled: one row inserted after "Docentes responsaveis:" (row 13) and one new
#    row appended at the end (row 22).
#  - Column layout: column A gets its own <col> band instead of being merged with column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the three existing formats (bold col A / wrapped col B / wrapped-red col C) via
# copy + paste-special-formats so no new (duplicate) style entries get created:
#   A:B:C triple, styles 1:2:3  -> template row 3
#   B:C pair,     styles   2:3  -> template row 2
#   A only,       style    1    -> template row 12

# --- Row 10 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Range("C10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Rows.Item(10).RowHeight = 60

# --- Row 11 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art"
$ws.Range("C11").Value = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art"
$ws.Rows.Item(11).RowHeight = 60

# --- Row 12 ---
$ws.Range("A12:A12").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Rows.Item(12).AutoFit()

# --- Row 13 ---
$ws.Range("B2:C2").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("A13").Clear()
$ws.Rows.Item(13).AutoFit()

# --- Row 14 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "A definir de acordo com o tópico programado"
$ws.Range("C14").Value = "A definir de acordo com o tópico programado"
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "To be defined according to the scheduled topic"
$ws.Range("C15").Value = "To be defined according to the scheduled topic"
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Range("C16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("C17").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18 ---
$ws.Range("A12:A12").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# --- Row 19 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."
$ws.Range("C22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."
$ws.Rows.Item(22).RowHeight = 120

$excel.CutCopyMode = $false

# Column layout fix: split the merged "A:B" <col> band (min=1,max=2) so column A gets its own
# min=1,max=1 entry (matches the diff). Re-stating column B's width is what forces the engine
# to split the band; column A's stored width (30.7109375) is left untouched.
$ws.Columns.Item(2).ColumnWidth = 60.7109375

$ws.Range("A1").Select()